$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the existing rows 18-24 (A:C) — these are the K_SERIES..K_ZUORDN entries that need
# to shift down by one row to make room for the new K_SEKTOREN row at 18. Read via Value2 and
# write with Value, without ever using a row Insert/shift (which reallocates an extra, unused
# style slot in this runtime) - this keeps styles.xml untouched, matching the target diff which
# only changes sheetData.
$snapshot = @()
for ($r = 18; $r -le 24; $r++) {
    $snapshot += ,@($ws.Cells.Item($r, 1).Value2, $ws.Cells.Item($r, 2).Value2, $ws.Cells.Item($r, 3).Value2)
}

for ($i = 0; $i -lt $snapshot.Length; $i++) {
    $destRow = 19 + $i
    $ws.Cells.Item($destRow, 1).Value = $snapshot[$i][0]
    $ws.Cells.Item($destRow, 2).Value = $snapshot[$i][1]
    $ws.Cells.Item($destRow, 3).Value = $snapshot[$i][2]
}

# Row 25 is brand new (the sheet previously ended at row 24), so it has no inherited formatting
# yet. Copy the formatting used by every other data row onto it.
$ws.Range("A17:C17").Copy()
$ws.Range("A25:C25").PasteSpecial(-4122)

# Fill the new row 18 with the K_SEKTOREN entry (row 18 already carries the correct data-row
# style, since that cell existed before and we only changed its value).
$ws.Cells.Item(18, 1).Value = "K_SEKTOREN"
$ws.Cells.Item(18, 2).Value = "Sektoren"
$ws.Cells.Item(18, 3).Value = "XXXSektoren"

$excel.CutCopyMode = $false
